# actualizacion base de datos en excel
# Update the "BD_actualizado" indicators table (Tabla33) on sheet
# "Indicadores mensuales" with the latest monthly data points for
# rows 248-254 (dates 2021-07 through 2023-01).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Indicadores mensuales")

$calcFormula = "=+Tabla33[[#This Row],[INFLATION]]+Tabla33[[#This Row],[DESEMPLEO]]-Tabla33[[#This Row],[PIB GROWTH]]"

# --- Row 248 (FECHA 2021-07-01): revise PIB GROWTH, ISM recalculates ---
$ws.Range("D248").Value = 0.1319

# --- Row 249 (FECHA 2021-08-01): revise PIB GROWTH and OIL PRICE ---
$ws.Range("D249").Value = 0.1319
$ws.Range("K249").Value = 67.71

# --- Row 250 (FECHA 2021-09-01): fill DESEMPLEO..SMMLV&TTE, revise OIL PRICE ---
$ws.Range("B250").Value = 0.1211
$ws.Range("C250").Value = 0.0451
$ws.Range("D250").Value = 0.1319
$ws.Range("E250").Formula = $calcFormula
$ws.Range("F250").Value = 3822
$ws.Range("G250").Value = 1014980
$ws.Range("K250").Value = 71.3

# --- Row 251 (FECHA 2021-10-01): fill full row of indicators ---
$ws.Range("B251").Value = 0.1179
$ws.Range("C251").Value = 0.0458
$ws.Range("D251").Value = 0.097
$ws.Range("D251").NumberFormat = "0.000"
$ws.Range("D251").HorizontalAlignment = -4108
$ws.Range("D251").Font.Color = 255
$ws.Range("D251").Font.Name = "Calibri"
$ws.Range("E251").Formula = $calcFormula
$ws.Range("F251").Value = 3773
$ws.Range("G251").Value = 1014980
$ws.Range("H251").Value = -0.013
$ws.Range("I251").Value = 0.106
$ws.Range("J251").Value = -0.192
$ws.Range("K251").Value = 81.22
$ws.Range("O251").Value = "Halloween"
$ws.Range("P251").Value = "Semana receso"
$ws.Range("Q251").Value = "Lluvias"
$ws.Range("R251").Value = "Iguales"
$ws.Range("S251").Value = "5% inflacion"
$ws.Range("T251").Value = "ICA. Renta personas"
$ws.Range("U251").Value = 60940
$ws.Range("V251").Value = 10038
$ws.Range("W251").Value = 10707

# --- Row 252 (FECHA 2021-11-01): fill full row of indicators ---
$ws.Range("B252").Value = 0.1083
$ws.Range("C252").Value = 0.0526
$ws.Range("D252").Value = 0.097
$ws.Range("D252").NumberFormat = "0.000"
$ws.Range("D252").HorizontalAlignment = -4108
$ws.Range("D252").Font.Color = 255
$ws.Range("D252").Font.Name = "Calibri"
$ws.Range("E252").Formula = $calcFormula
$ws.Range("F252").Value = 3903
$ws.Range("G252").Value = 1014980
$ws.Range("H252").Value = -0.014
$ws.Range("I252").Value = 0.064
$ws.Range("J252").Value = -0.131
$ws.Range("K252").Value = 78.6
$ws.Range("O252").Value = "Navidad - Black Friday"
$ws.Range("P252").Value = "No"
$ws.Range("Q252").Value = "Lluvias"
$ws.Range("R252").Value = "Iguales"
$ws.Range("S252").Value = "5% inflacion"
$ws.Range("T252").Value = "ICA"
$ws.Range("U252").Value = 60321
$ws.Range("V252").Value = 8066
$ws.Range("W252").Value = 9173

# --- Row 253 (FECHA 2021-12-01): fill full row of indicators ---
$ws.Range("B253").Value = 0.1101
$ws.Range("C253").Value = 0.0562
$ws.Range("D253").Value = 0.097
$ws.Range("D253").NumberFormat = "0.000"
$ws.Range("D253").HorizontalAlignment = -4108
$ws.Range("D253").Font.Color = 255
$ws.Range("D253").Font.Name = "Calibri"
$ws.Range("E253").Formula = $calcFormula
$ws.Range("F253").Value = 3963
$ws.Range("G253").Value = 1014980
$ws.Range("H253").Value = -0.07
$ws.Range("I253").Value = 0.021
$ws.Range("J253").Value = -0.206
$ws.Range("K253").Value = 71.69
$ws.Range("O253").Value = "Navidad"
$ws.Range("P253").Value = "A. B y universidades"
$ws.Range("Q253").Value = "Soleado"
$ws.Range("R253").Value = "Prima"
$ws.Range("S253").Value = "5% inflacion"
$ws.Range("T253").Value = "ICA"
$ws.Range("U253").Value = 77651
$ws.Range("V253").Value = 11643
$ws.Range("W253").Value = 9143

# --- Row 254 (FECHA 2022-01-01): fill TRM, SMMLV&TTE and ICC/IEE/ICE ---
$ws.Range("F254").Value = 4000
$ws.Range("G254").Value = 1117172
$ws.Range("H254").Value = -0.135
$ws.Range("I254").Value = -0.024
$ws.Range("J254").Value = -0.3
